$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'29.197.53"
$ws.Range("E2").Value = "  -0.63%  "

# Row 3
$ws.Range("D3").Value = "'1.860.95"
$ws.Range("E3").Value = "  -0.97%  "

# Row 4
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.08%  "

# Row 5
$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").Value = "'241.92"
$ws.Range("E5").Value = "  -0.75%  "

# Row 6
$ws.Range("B6").Value = "XRP"
$ws.Range("C6").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D6").Value = "'0.7019"
$ws.Range("E6").Value = "  -2.05%  "

# Row 7
$ws.Range("D7").Value = "'1.001"
$ws.Range("E7").Value = "  -0.03%  "

# Row 8
$ws.Range("E8").Value = "  -1.96%  "

# Row 9
$ws.Range("E9").Value = "  -1.15%  "

# Row 10
$ws.Range("D10").Value = "'23.85"
$ws.Range("E10").Value = "  -4.24%  "

# Row 11
$ws.Range("D11").Value = "'0.07803"
$ws.Range("E11").Value = "  -3.45%  "

# Row 12
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "'1.835.65"
$ws.Range("E12").Value = "  -2.26%  "

# Row 13
$ws.Range("B13").Value = "Litecoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D13").Value = "'92.71"
$ws.Range("E13").Value = "  -2.15%  "

# Row 14
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'5.123"
$ws.Range("E14").Value = "  -1.95%  "

# Row 15
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").Value = "'0.6909"
$ws.Range("E15").Value = "  -2.44%  "

# Row 16
$ws.Range("D16").Value = "'6.539"
$ws.Range("E16").Value = "  +2.37%  "

# Row 17
$ws.Range("D17").Value = "'0.000008439"
$ws.Range("E17").Value = "  +0.11%  "

# Row 18
$ws.Range("D18").Value = "'29.260.13"
$ws.Range("E18").Value = "  -0.42%  "

# Row 19
$ws.Range("D19").Value = "'250.29"
$ws.Range("E19").Value = "  -0.95%  "

# Row 20
$ws.Range("D20").Value = "'2.116.20"
$ws.Range("E20").Value = "  -0.39%  "

# Row 21
$ws.Range("D21").Value = "'12.93"
$ws.Range("E21").Value = "  -3.17%  "

# Row 22
$ws.Range("D22").Value = "'0.9998"
$ws.Range("E22").Value = "  -0.11%  "

# Row 23
$ws.Range("D23").Value = "'7.595"
$ws.Range("E23").Value = "  -1.09%  "

# Row 24
$ws.Range("D24").Value = "'1.000"
$ws.Range("E24").Value = "  -0.02%  "

# Row 25
$ws.Range("D25").Value = "'0.1531"
$ws.Range("E25").Value = "  -2.97%  "

# Row 26
$ws.Range("D26").Value = "'160.23"
$ws.Range("E26").Value = "  -1.04%  "

# Row 27
$ws.Range("E27").Value = "  -1.96%  "

# Row 28
$ws.Range("D28").Value = "'18.58"
$ws.Range("E28").Value = "  -2.20%  "

# Row 29
$ws.Range("D29").Value = "'1.568"
$ws.Range("E29").Value = "  +3.92%  "

# Row 30
$ws.Range("D30").Value = "'4.278"
$ws.Range("E30").Value = "  -3.18%  "

# Row 31
$ws.Range("E31").Value = "  -1.54%  "

# Row 32
$ws.Range("E32").Value = "  -1.19%  "

# Row 33
$ws.Range("D33").Value = "'0.05224"
$ws.Range("E33").Value = "  -1.59%  "

# Row 34
$ws.Range("D34").Value = "'1.874"
$ws.Range("E34").Value = "  -3.55%  "

# Row 35
$ws.Range("D35").Value = "'0.7566"
$ws.Range("E35").Value = "  -0.03%  "

# Row 36
$ws.Range("E36").Value = "  +0.11%  "

# Row 37
$ws.Range("D37").Value = "'2.708"
$ws.Range("E37").Value = "  +0.20%  "

# Row 38
$ws.Range("D38").Value = "'0.01861"
$ws.Range("E38").Value = "  -1.29%  "

# Row 39
$ws.Range("D39").Value = "'1.223.13"
$ws.Range("E39").Value = "  -4.66%  "

# Row 40
$ws.Range("D40").Value = "'2.721"
$ws.Range("E40").Value = "  -1.57%  "

# Row 41
$ws.Range("D41").Value = "'0.9020"
$ws.Range("E41").Value = "  -0.71%  "

# Row 42
$ws.Range("D42").Value = "'110.11"
$ws.Range("E42").Value = "  -1.36%  "

# Row 43
$ws.Range("D43").Value = "'5.791"
$ws.Range("E43").Value = "  -9.59%  "

# Row 44
$ws.Range("D44").Value = "'1.0000"
$ws.Range("E44").Value = "  -0.07%  "

# Row 45
$ws.Range("D45").Value = "'2.012.60"
$ws.Range("E45").Value = "  -0.73%  "

# Row 46
$ws.Range("E46").Value = "  -3.87%  "

# Row 47
$ws.Range("D47").Value = "'65.03"
$ws.Range("E47").Value = "  -12.27%  "

# Row 48
$ws.Range("D48").Value = "'0.5189"
$ws.Range("E48").Value = "  -0.36%  "

# Row 49
$ws.Range("D49").Value = "'9.508"
$ws.Range("E49").Value = "  -0.27%  "

# Row 50
$ws.Range("D50").Value = "'1.768"
$ws.Range("E50").Value = "  -2.14%  "

# Row 51
$ws.Range("D51").Value = "'7.032"
$ws.Range("E51").Value = "  -0.95%  "
